{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Goal (per the commit's neutral-language update):\n//   \"affecting all Black and Asian-American voters\" -> \"affecting 50M voters\"\n// in three places:\n//   1. PROFESSIONAL SUMMARY paragraph               (plain-text run)\n//   2. Siege Analytics bullet point                 (\"50M\" becomes a bold,\n//      colored run matching the existing \"23%\"/\"64%\" runs in that bullet)\n//   3. \"Impact: Corrected demographic data ...\" line (plain-text run,\n//      also appends \"nationwide\" after \"50M voters\")\n//\n// NOTE: Change #2 is applied FIRST, while \"affecting all Black and\n// Asian-American voters\" is still unique to that bullet's search text and\n// \"50M\" does not yet appear anywhere else in the document. That lets us\n// uniquely locate the freshly-inserted \"50M\" run and bold/color it before\n// changes #1 and #3 introduce their own plain-text \"50M\" occurrences.\n\nconst OLD_PHRASE = \"affecting all Black and Asian-American voters\";\n\n// ---------------------------------------------------------------------\n// Change 2: bullet point under \"Partner - Siege Analytics\" \u2014 split into\n// three runs so \"50M\" can carry bold + the same accent color used by the\n// neighboring \"23%\"/\"64%\" statistic runs in that same bullet.\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n  const results = body.search(\n    OLD_PHRASE + \", developed geospatial machine learning\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for the Siege Analytics bullet, found \" +\n        results.items.length\n    );\n  }\n\n  // Replace the whole matched span with the new plain text first ...\n  results.items[0].insertText(\n    \"affecting 50M voters, developed geospatial machine learning\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  // ... then locate the \"50M\" that was just inserted (unique at this point)\n  // and apply the bold + accent-color formatting used elsewhere in the doc.\n  const body2 = context.document.body;\n  const fiftyM = body2.search(\"50M\", { matchCase: true });\n  fiftyM.load(\"items\");\n  await context.sync();\n\n  if (fiftyM.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '50M' while formatting, found \" +\n        fiftyM.items.length\n    );\n  }\n\n  fiftyM.items[0].font.bold = true;\n  fiftyM.items[0].font.color = \"2C3E50\";\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 1: PROFESSIONAL SUMMARY paragraph \u2014 plain text replace, no\n// formatting change.\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n  const results = body.search(\n    OLD_PHRASE + \", developed geospatial ML\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for the PROFESSIONAL SUMMARY sentence, found \" +\n        results.items.length\n    );\n  }\n\n  results.items[0].insertText(\n    \"affecting 50M voters, developed geospatial ML\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 3: \"Impact: Corrected demographic data ...\" line \u2014 plain text\n// replace, no formatting change. Also appends \"nationwide\".\n// ---------------------------------------------------------------------\n{\n  const body = context.document.body;\n  const results = body.search(\n    \"Impact: Corrected demographic data \" + OLD_PHRASE + \", improved\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for the Impact line, found \" +\n        results.items.length\n    );\n  }\n\n  results.items[0].insertText(\n    \"Impact: Corrected demographic data affecting 50M voters nationwide, improved\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Goal (per the commit's neutral-language update):\n#   \"affecting all Black and Asian-American voters\" -> \"affecting 50M voters\"\n# in three places:\n#   1. PROFESSIONAL SUMMARY paragraph               (plain-text run)\n#   2. Siege Analytics bullet point                 (\"50M\" becomes a bold,\n#      colored run matching the existing \"23%\"/\"64%\" runs in that bullet)\n#   3. \"Impact: Corrected demographic data ...\" line (plain-text run,\n#      also appends \"nationwide\" after \"50M voters\")\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# Helper: find a unique, literal phrase anywhere in the document and\n# replace it outright (Find/Replace, wdReplaceAll scoped to one hit by\n# the phrase's uniqueness). No formatting is touched.\n# -----------------------------------------------------------------------\nfunction Replace-UniqueText {\n    param($doc, [string]$oldText, [string]$newText)\n\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0                 # wdFindStop - do not wrap around\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $ok = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) # wdReplaceAll\n    if (-not $ok) {\n        throw \"Replace-UniqueText: phrase not found -> $oldText\"\n    }\n}\n\n# -----------------------------------------------------------------------\n# Change 2: bullet point under \"Partner - Siege Analytics\" bullet list.\n# This one needs THREE runs in the output (\"50M\" bold + colored 2C3E50,\n# matching the neighboring \"23%\"/\"64%\" statistic runs in the same\n# bullet), so it can't be done as a plain Find/Replace.\n# -----------------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"affecting all Black and Asian-American voters, developed geospatial machine learning\"\n$find2.Forward = $true\n$find2.Wrap = 0\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n\n$found2 = $find2.Execute()\nif (-not $found2) {\n    throw \"Change 2: Siege Analytics bullet phrase not found\"\n}\n\n$bulletRange = $find2.Parent\n$bulletStart = $bulletRange.Start\n\n# Overwrite the whole matched span with the new plain text (this shrinks\n# the range; re-derive the end from the range itself afterwards).\n$bulletRange.Text = \"affecting 50M voters, developed geospatial machine learning\"\n$bulletEnd = $bulletRange.End\n\n# Re-fetch a fresh Range over the same span (the just-written $bulletRange\n# object's cached .Text can read stale) and scope a Find to it so \"50M\" is\n# located precisely without hand-computed character offsets.\n$freshRange = $d.Range($bulletStart, $bulletEnd)\n$fiftyMFind = $freshRange.Find\n$fiftyMFind.ClearFormatting()\n$fiftyMFind.Text = \"50M\"\n$fiftyMFind.Forward = $true\n$fiftyMFind.Wrap = 0\n$fiftyMFind.Format = $false\n$fiftyMFind.MatchCase = $true\n\n$foundFiftyM = $fiftyMFind.Execute()\nif (-not $foundFiftyM) {\n    throw \"Change 2: could not re-locate '50M' after replace\"\n}\n\n$fiftyMRange = $fiftyMFind.Parent\n$fiftyMRange.Font.Bold = $true\n$fiftyMRange.Font.Color = \"2C3E50\"   # matches the other bolded stat runs\n\n# -----------------------------------------------------------------------\n# Change 1: PROFESSIONAL SUMMARY paragraph - plain text replace.\n# -----------------------------------------------------------------------\nReplace-UniqueText $d `\n    \"affecting all Black and Asian-American voters, developed geospatial ML\" `\n    \"affecting 50M voters, developed geospatial ML\"\n\n# -----------------------------------------------------------------------\n# Change 3: \"Impact: Corrected demographic data ...\" line - plain text\n# replace (also adds \"nationwide\").\n# -----------------------------------------------------------------------\nReplace-UniqueText $d `\n    \"Impact: Corrected demographic data affecting all Black and Asian-American voters, improved\" `\n    \"Impact: Corrected demographic data affecting 50M voters nationwide, improved\"\n"}
